$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Arial 11pt font used for the cleaned-up "sample number" column (G)
$sampleFont = "E7420"

for ($r = 2; $r -le 27; $r++) {
    $gCell = $ws.Cells.Item($r, 7)
    $gCell.Value = $sampleFont
    $gCell.Font.Name = "Arial"
    $gCell.Font.Size = 11
    $gCell.Font.Color = 0

    $hCell = $ws.Cells.Item($r, 8)
    $hCell.Formula = "=FALSE()"
}

$ws.Range("G2:G27").Select()
